$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

$ws.Range("D2").Value = "37.890.55"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "2.103.92"
$ws.Range("E3").Value = "  +1.84%  "
$ws.Range("E4").Value = "  +0.05%  "
Set-TextValue "D5" "235.60"
$ws.Range("E5").Value = "  +0.94%  "
Set-TextValue "D6" "0.624"
$ws.Range("E6").Value = "  +0.93%  "
Set-TextValue "D7" "58.45"
$ws.Range("E7").Value = "  -1.78%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +1.96%  "
Set-TextValue "D10" "0.0782"
$ws.Range("E10").Value = "  +2.52%  "
$ws.Range("E11").Value = "  +1.25%  "
$ws.Range("D12").Value = "2.412.36"
$ws.Range("E12").Value = "  +1.85%  "
Set-TextValue "D13" "14.46"
$ws.Range("E13").Value = "  +0.38%  "
Set-TextValue "D14" "21.16"
$ws.Range("E14").Value = "  -0.04%  "
Set-TextValue "D15" "0.785"
$ws.Range("E15").Value = "  +1.07%  "
Set-TextValue "D16" "5.24"
$ws.Range("E16").Value = "  +1.22%  "
$ws.Range("D17").Value = "2.098.99"
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("D18").Value = "37.818.17"
$ws.Range("E18").Value = "  +1.13%  "
Set-TextValue "D19" "6.28"
$ws.Range("E19").Value = "  +0.74%  "
Set-TextValue "D20" "70.29"
$ws.Range("E20").Value = "  +1.59%  "
$ws.Range("D21").Value = "0.0₃0823"
$ws.Range("E21").Value = "  +1.27%  "
Set-TextValue "D22" "227.56"
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  +0.74%  "
Set-TextValue "D25" "2.42"
$ws.Range("E25").Value = "  +0.06%  "
Set-TextValue "D26" "168.07"
$ws.Range("E26").Value = "  +1.33%  "
Set-TextValue "D27" "8.99"
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("E28").Value = "  +3.39%  "
$ws.Range("E29").Value = "  -4.39%  "
$ws.Range("E30").Value = "  +1.89%  "
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("E32").Value = "  +3.61%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D33" "0.0624"
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D34" "2.58"
$ws.Range("E34").Value = "  -0.53%  "
Set-TextValue "D35" "4.57"
$ws.Range("E35").Value = "  +0.02%  "
Set-TextValue "D36" "3.46"
$ws.Range("E36").Value = "  +4.65%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("E38").Value = "  +0.14%  "
Set-TextValue "D39" "5.51"
$ws.Range("E39").Value = "  -7.23%  "
$ws.Range("E40").Value = "  +1.89%  "
$ws.Range("E41").Value = "  -0.37%  "
Set-TextValue "D42" "97.07"
$ws.Range("E42").Value = "  +0.44%  "
$ws.Range("D43").Value = "1.470.23"
$ws.Range("E43").Value = "  +0.42%  "
Set-TextValue "D44" "0.0213"
$ws.Range("E44").Value = "  +0.72%  "
Set-TextValue "D45" "1.17"
$ws.Range("E45").Value = "  +0.07%  "
Set-TextValue "D46" "4.19"
$ws.Range("E46").Value = "  -10.89%  "
$ws.Range("E47").Value = "  +2.12%  "
Set-TextValue "D48" "15.46"
$ws.Range("E48").Value = "  -1.71%  "
Set-TextValue "D49" "7.33"
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("E50").Value = "  +2.84%  "
$ws.Range("D51").Value = "2.299.15"
$ws.Range("E51").Value = "  +1.93%  "

"done"